$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 15 (pushes the existing rows 15-53 down to 16-54,
# and the worksheet dimension grows from A1:T53 to A1:T54).
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A15").Value = 5
$ws.Range("B15").Value = "Macroferia Regional de Talca"
$ws.Range("C15").Value = "Maule"
$ws.Range("D15").Value = 44581
$ws.Range("E15").Value = 7
$ws.Range("F15").Value = "Fruta"
$ws.Range("G15").Value = 100101
$ws.Range("H15").Value = "Berries"
$ws.Range("I15").Value = 100101001
$ws.Range("J15").Value = "Arándano (blue)"
$ws.Range("K15").Value = "Sin especificar"
$ws.Range("L15").Value = "Segunda"
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 3400
$ws.Range("O15").Value = 3400
$ws.Range("P15").Value = 3400
$ws.Range("Q15").Value = "$/bandeja 2 kilos"
$ws.Range("R15").Value = "Provincia de Linares"
$ws.Range("S15").Value = 1700
$ws.Range("T15").Value = 2
